$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final table data (header in row 1 stays the same).
# The "Nick Richards" row has been removed and the remaining rows
# re-ordered to their final positions.
$data = @(
    @("Tyrese Haliburton", "PG,SG", "Indiana Pacers"),
    @("Klay Thompson", "SG,SF", "Dallas Mavericks"),
    @("Stephen Curry", "PG,SG", "Golden State Warriors"),
    @("Austin Reaves", "PG,SG", "Los Angeles Lakers"),
    @("OG Anunoby", "SF,PF", "New York Knicks"),
    @("Franz Wagner", "SF,PF", "Orlando Magic"),
    @("Kevin Durant", "SF,PF", "Phoenix Suns"),
    @("Cameron Johnson", "SF,PF", "Brooklyn Nets"),
    @("Karl-Anthony Towns", "PF,C", "New York Knicks"),
    @("Jarrett Allen", "C", "Cleveland Cavaliers"),
    @("Trey Murphy III", "SF,PF", "New Orleans Pelicans"),
    @("Jalen Duren", "C", "Detroit Pistons"),
    @("Keegan Murray", "SF,PF", "Sacramento Kings"),
    @("Mark Williams", "C", "Charlotte Hornets"),
    @("Tyrese Maxey", "PG,SG", "Philadelphia 76ers"),
    @("Darius Garland", "PG", "Cleveland Cavaliers")
)

# The old sheet had 17 data rows (rows 2-18); the new sheet has 16 (rows 2-17).
# Clear the old data range first so the extra trailing row is removed.
$ws.Range("A2:C18").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
}
